$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Entry 6 (row 6): update compound equivalents and drop the "(1:1)" suffix
# from the amine name.
$ws.Range("B6").Value = 1.6
$ws.Range("C6").Value = "Benzenamine hydrochloride"
$ws.Range("F6").Value = 2.5
$ws.Range("H6").Value = 3.3

# Move the active selection to A6, matching the saved view state.
$ws.Range("A6").Select()
